$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 21:52"

# Row 4
$ws.Range("B4").Value = 2970468
$ws.Range("C4").Value = 34698
$ws.Range("D4").Value = 1280758
$ws.Range("E4").Value = 1557200
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 192
$ws.Range("H4").Value = 132510

# Row 6
$ws.Range("B6").Value = 697836
$ws.Range("C6").Value = 23932
$ws.Range("D6").Value = 424891
$ws.Range("E6").Value = 253245
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 421
$ws.Range("H6").Value = 19700

# Row 18
$ws.Range("B18").Value = 197558
$ws.Range("C18").Value = 140
$ws.Range("D18").Value = 181700
$ws.Range("E18").Value = 6772
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 9086

# Row 19
$ws.Range("B19").Value = 196750
$ws.Range("C19").Value = 8773
$ws.Range("D19").Value = 93315
$ws.Range("E19").Value = 100236
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 173
$ws.Range("H19").Value = 3199

# Row 90
$ws.Range("A90").Value = "Costa Rica"
$ws.Range("B90").Value = 4996
$ws.Range("C90").Value = 375
$ws.Range("D90").Value = 1745
$ws.Range("E90").Value = 3232
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 19

# Row 91
$ws.Range("A91").Value = "Bosnia y Herzegovina"
$ws.Range("B91").Value = 4962
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 2550
$ws.Range("E91").Value = 2221
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 191

# Row 92
$ws.Range("A92").Value = "Guayana Francesa"
$ws.Range("B92").Value = 4913
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 1866
$ws.Range("E92").Value = 3031
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 16

# Row 93
$ws.Range("A93").Value = "Mauritania"
$ws.Range("B93").Value = 4827
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 1805
$ws.Range("E93").Value = 2893
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 129

# Row 94
$ws.Range("A94").Value = "Republica de Yibuti"
$ws.Range("B94").Value = 4792
$ws.Range("C94").Value = 56
$ws.Range("D94").Value = 4593
$ws.Range("E94").Value = 144
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 55

# Row 96
$ws.Range("A96").Value = "Estado de Palestina"
$ws.Range("B96").Value = 4277
$ws.Range("C96").Value = 442
$ws.Range("D96").Value = 491
$ws.Range("E96").Value = 3770
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 3
$ws.Range("H96").Value = 16

# Row 97
$ws.Range("A97").Value = "Hungria"
$ws.Range("B97").Value = 4183
$ws.Range("C97").Value = 9
$ws.Range("D97").Value = 2811
$ws.Range("E97").Value = 783
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 589

# Row 130
$ws.Range("D130").Value = 1048
$ws.Range("E130").Value = 90

# Row 205
$ws.Range("A205").Value = "Dominica"

# Row 206
$ws.Range("A206").Value = "Fiyi"
